$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Step 1: Normalize formatting (styles) before overwriting values ----
# IMPORTANT: capture the "style 2" template (G9, blank) for G15 FIRST,
# before G9 formatting gets changed later in this script.
$ws.Range("G9").Copy()
$ws.Range("G15").PasteSpecial(-4122)

# Style "1" template (A9:F9, already style 1) -> apply to A8:F8 (was style 2)
$ws.Range("A9:F9").Copy()
$ws.Range("A8:F8").PasteSpecial(-4122)

# Rows 10-15 A:F need style "1" -> copy from A9:F9 template
$ws.Range("A9:F9").Copy()
$ws.Range("A10:F10").PasteSpecial(-4122)
$ws.Range("A9:F9").Copy()
$ws.Range("A11:F11").PasteSpecial(-4122)
$ws.Range("A9:F9").Copy()
$ws.Range("A12:F12").PasteSpecial(-4122)
$ws.Range("A9:F9").Copy()
$ws.Range("A13:F13").PasteSpecial(-4122)
$ws.Range("A9:F9").Copy()
$ws.Range("A14:F14").PasteSpecial(-4122)
$ws.Range("A9:F9").Copy()
$ws.Range("A15:F15").PasteSpecial(-4122)

# Style "1" template (A9) -> apply to G14 (needs style 1)
$ws.Range("A9").Copy()
$ws.Range("G14").PasteSpecial(-4122)

# Style "0" (default/no special style) template -> A2:G2 (already style 0)
# Apply to rows 5,6,7 (remove style 1/2 -> default)
$ws.Range("A2:G2").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)
$ws.Range("A6:G6").PasteSpecial(-4122)
$ws.Range("A7:G7").PasteSpecial(-4122)

# Style "0" template (A2, default) -> apply to G8:G13 (remove style 2 -> default)
$ws.Range("A2").Copy()
$ws.Range("G8").PasteSpecial(-4122)
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("G10").PasteSpecial(-4122)
$ws.Range("G11").PasteSpecial(-4122)
$ws.Range("G12").PasteSpecial(-4122)
$ws.Range("G13").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---- Step 2: Set cell values (this also governs shared-string append order) ----
# Row 2
$ws.Cells.Item(2,1).Value = "math"
$ws.Cells.Item(2,2).Value = 2025
$ws.Cells.Item(2,3).Value = "M"
$ws.Cells.Item(2,4).Value = "qp-202505-mathematics-p43"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = "subjects/math/2025/M/qp-202505-mathematics-p43"
$ws.Cells.Item(2,7).Value = "1.html"

# Row 3
$ws.Cells.Item(3,1).Value = "math"
$ws.Cells.Item(3,2).Value = 2025
$ws.Cells.Item(3,3).Value = "M"
$ws.Cells.Item(3,4).Value = "qp-202505-mathematics-p43"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = "subjects/math/2025/M/qp-202505-mathematics-p43"
$ws.Cells.Item(3,7).Value = "2.html"

# Row 4
$ws.Cells.Item(4,1).Value = "math"
$ws.Cells.Item(4,2).Value = 2025
$ws.Cells.Item(4,3).Value = "M"
$ws.Cells.Item(4,4).Value = "qp-202505-mathematics-p43"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = "subjects/math/2025/M/qp-202505-mathematics-p43"
$ws.Cells.Item(4,7).Value = "3.html"

# Row 5
$ws.Cells.Item(5,1).Value = "math"
$ws.Cells.Item(5,2).Value = 2025
$ws.Cells.Item(5,3).Value = "M"
$ws.Cells.Item(5,4).Value = "qp-202505-mathematics-p43"
$ws.Cells.Item(5,5).Value = 4
$ws.Cells.Item(5,6).Value = "subjects/math/2025/M/qp-202505-mathematics-p43"
$ws.Cells.Item(5,7).Value = "4.html"

# Row 6
$ws.Cells.Item(6,1).Value = "math"
$ws.Cells.Item(6,2).Value = 2025
$ws.Cells.Item(6,3).Value = "M"
$ws.Cells.Item(6,4).Value = "qp-202505-mathematics-p43"
$ws.Cells.Item(6,5).Value = 5
$ws.Cells.Item(6,6).Value = "subjects/math/2025/M/qp-202505-mathematics-p43"
$ws.Cells.Item(6,7).Value = "5.html"

# Row 7
$ws.Cells.Item(7,1).Value = "math"
$ws.Cells.Item(7,2).Value = 2025
$ws.Cells.Item(7,3).Value = "M"
$ws.Cells.Item(7,4).Value = "qp-202505-mathematics-p43"
$ws.Cells.Item(7,5).Value = 6
$ws.Cells.Item(7,6).Value = "subjects/math/2025/M/qp-202505-mathematics-p43"
$ws.Cells.Item(7,7).Value = "6.html"

# Row 8
$ws.Cells.Item(8,1).Value = "math"
$ws.Cells.Item(8,2).Value = 2025
$ws.Cells.Item(8,3).Value = "M"
$ws.Cells.Item(8,4).Value = "qp-202505-mathematics-p45"
$ws.Cells.Item(8,5).Value = 1
$ws.Cells.Item(8,6).Value = "subjects/math/2025/M/qp-202505-mathematics-p45"
$ws.Cells.Item(8,7).Value = "1.html"

# Row 9
$ws.Cells.Item(9,1).Value = "math"
$ws.Cells.Item(9,2).Value = 2025
$ws.Cells.Item(9,3).Value = "M"
$ws.Cells.Item(9,4).Value = "qp-202505-mathematics-p45"
$ws.Cells.Item(9,5).Value = 2
$ws.Cells.Item(9,6).Value = "subjects/math/2025/M/qp-202505-mathematics-p45"
$ws.Cells.Item(9,7).Value = "2.html"

# Row 10
$ws.Cells.Item(10,1).Value = "math"
$ws.Cells.Item(10,2).Value = 2025
$ws.Cells.Item(10,3).Value = "M"
$ws.Cells.Item(10,4).Value = "qp-202505-mathematics-p45"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = "subjects/math/2025/M/qp-202505-mathematics-p45"
$ws.Cells.Item(10,7).Value = "3.html"

# Row 11
$ws.Cells.Item(11,1).Value = "math"
$ws.Cells.Item(11,2).Value = 2025
$ws.Cells.Item(11,3).Value = "M"
$ws.Cells.Item(11,4).Value = "qp-202505-mathematics-p45"
$ws.Cells.Item(11,5).Value = 4
$ws.Cells.Item(11,6).Value = "subjects/math/2025/M/qp-202505-mathematics-p45"
$ws.Cells.Item(11,7).Value = "4.html"

# Row 12
$ws.Cells.Item(12,1).Value = "math"
$ws.Cells.Item(12,2).Value = 2025
$ws.Cells.Item(12,3).Value = "M"
$ws.Cells.Item(12,4).Value = "qp-202505-mathematics-p45"
$ws.Cells.Item(12,5).Value = 5
$ws.Cells.Item(12,6).Value = "subjects/math/2025/M/qp-202505-mathematics-p45"
$ws.Cells.Item(12,7).Value = "5.html"

# Row 13
$ws.Cells.Item(13,1).Value = "math"
$ws.Cells.Item(13,2).Value = 2025
$ws.Cells.Item(13,3).Value = "M"
$ws.Cells.Item(13,4).Value = "qp-202505-mathematics-p45"
$ws.Cells.Item(13,5).Value = 6
$ws.Cells.Item(13,6).Value = "subjects/math/2025/M/qp-202505-mathematics-p45"
$ws.Cells.Item(13,7).Value = "6.html"

# Row 14
$ws.Cells.Item(14,1).Value = "math"
$ws.Cells.Item(14,2).Value = 2025
$ws.Cells.Item(14,3).Value = "M"
$ws.Cells.Item(14,4).Value = "qp-202505-mathematics-p45"
$ws.Cells.Item(14,5).Value = 7
$ws.Cells.Item(14,6).Value = "subjects/math/2025/M/qp-202505-mathematics-p45"
$ws.Cells.Item(14,7).Value = "7.html"

# Row 15
$ws.Cells.Item(15,1).Value = "math"
$ws.Cells.Item(15,2).Value = 2025
$ws.Cells.Item(15,3).Value = "S1"
$ws.Cells.Item(15,4).Value = "qp-202505-mathematics-p51"
$ws.Cells.Item(15,5).Value = 1
$ws.Cells.Item(15,6).Value = "subjects/math/2025/M/qp-202505-mathematics-p45"
$ws.Cells.Item(15,7).ClearContents()

# ---- Step 3: Remove stray H8 ("y") cell content ----
$ws.Cells.Item(8,8).ClearContents()

# ---- Step 4: Update selection to match final cursor position ----
$ws.Range("F21").Select()

Write-Host "Edit complete"